$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The recorded interview Q&A/evaluation details for user1 (columns D:M of row 2 -
# test_taken, final_rating, answer_1, evaluation_1, answer_2, evaluation_2,
# answer_3, evaluation_3, answer_4, evaluation_4) are no longer wanted, so wipe
# them out while leaving the username/interview_type columns (A2:B2) untouched.
$ws.Range("D2:M2").ClearContents()

# Leave the cursor/selection parked on G11, as in the final saved state.
$ws.Range("G11").Select() | Out-Null
